# Commit: "calling it a platform"
#
# 1) Bump the cached "datetimeFigureOut" footer-date text from 14.10.2024
#    to 15.10.2024 everywhere it is cached (the slide master and every
#    slide layout).
# 2) On slide 1, add a second line "                      Platform" under
#    the existing "                 Infrastructure" line in the
#    "Rectangle 3" shape, matching its formatting.

$p = $ppt.ActivePresentation

# --- 1) Update the cached date field text (master + all layouts) -------
$oldDate = "14.10.2024"
$newDate = "15.10.2024"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 2) Add the "Platform" line under "Infrastructure" on slide 1 ------
$slide1 = $p.Slides.Item(1)
$rect3 = $slide1.Shapes.Item(1)
[void]$rect3.TextFrame.TextRange.InsertAfter("`r                      Platform")
